$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.002847671508789
$ws.Range("B1").Value = 4.287489891052246
$ws.Range("C1").Value = 3.759507179260254
$ws.Range("D1").Value = 1.776339888572693
$ws.Range("E1").Value = 0.7869437336921692
